# Changed trade_id logic for exit
#
# Adds 5 new "MP41" trade rows (entry == exit, zero P&L, a flat tax-only
# net loss of -35.40) to the MPWizard trade log, and mirrors each one into
# the DTD running-balance ledger. The last of the five rows records its
# entry/exit time as a literal text timestamp instead of a date serial
# (reflects the trade_id/exit-time logic change called out in the commit).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# MPWizard: duplicate the formatting of the last existing row (row 11)
# down into rows 12-16, then overwrite with the new trade data.
# ---------------------------------------------------------------------
$mp = $wb.Worksheets.Item("MPWizard")

for ($r = 12; $r -le 16; $r++) {
    $mp.Range("A11:N11").Copy($mp.Range("A" + $r + ":N" + $r))
}

$entryExit = 45235.59791666667

for ($r = 12; $r -le 16; $r++) {
    $mp.Cells.Item($r, 1).Value = "MP41"
    $mp.Cells.Item($r, 2).Value = "NIFTY09NOV23P19250"
    $mp.Cells.Item($r, 3).Value = "Long"
    $mp.Cells.Item($r, 6).Value = 0
    $mp.Cells.Item($r, 7).Value = 0
    $mp.Cells.Item($r, 8).Value = 0
    $mp.Cells.Item($r, 9).Value = 0
    $mp.Cells.Item($r, 10).Value = 0
    $mp.Cells.Item($r, 11).Value = 50
    $mp.Cells.Item($r, 12).Value = 0
    $mp.Cells.Item($r, 13).Value = 35.4
    $mp.Cells.Item($r, 14).Value = -35.4
}

# Rows 12-15: entry/exit time as a real date-time serial (same instant)
for ($r = 12; $r -le 15; $r++) {
    $mp.Cells.Item($r, 4).Value = $entryExit
    $mp.Cells.Item($r, 5).Value = $entryExit
}

# Row 16: entry/exit time recorded as literal text instead of a serial.
# Copy column A's (plain centered) format onto D16:E16 before writing the
# text so it doesn't keep the D/E date-time format, then fill in the value.
$mp.Range("A16").Copy()
$mp.Range("D16:E16").PasteSpecial(-4122)
$mp.Cells.Item(16, 4).Value = "2023-11-05 14:21:00"
$mp.Cells.Item(16, 5).Value = "2023-11-05 14:21:00"

# ---------------------------------------------------------------------
# DTD: mirror the 5 new trades into the running-balance ledger,
# each one a tax-only loss of ₹35.40 against the prior balance.
# ---------------------------------------------------------------------
$dtd = $wb.Worksheets.Item("DTD")

for ($r = 15; $r -le 19; $r++) {
    $dtd.Range("A14:G14").Copy($dtd.Range("A" + $r + ":G" + $r))
}

$balances = @(18820.62, 18785.22, 18749.82, 18714.42, 18679.02)
$slNos = @(14, 15, 16, 17, 18)

for ($i = 0; $i -lt 5; $i++) {
    $r = 15 + $i

    $dtd.Cells.Item($r, 2).NumberFormat = "@"
    $dtd.Cells.Item($r, 3).NumberFormat = "@"

    $dtd.Cells.Item($r, 1).Value = $slNos[$i]
    $dtd.Cells.Item($r, 2).Value = "05-Nov-23"
    $dtd.Cells.Item($r, 3).Value = "Sunday"
    $dtd.Cells.Item($r, 4).Value = "MP41"
    $dtd.Cells.Item($r, 5).Value = "MPWizard"
    $dtd.Cells.Item($r, 6).Value = "-" + [char]8377 + "35.40"
    $dtd.Cells.Item($r, 7).Value = [char]8377 + ([string]::Format("{0:N2}", $balances[$i]))

    $dtd.Range("A" + $r + ":G" + $r).ClearFormats()
}

# ---------------------------------------------------------------------
# Workbook-level housekeeping: reset active tab back to the first sheet.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Stocks").Activate()
